# "template download according to add,delete"
#
# Adds an "age" column, splits the old combined legend column into a
# proper legend + Date of joining + Date of Leaving set of headers, adds
# a new employee row (Abdul Quadir) above the existing Avinash row, and
# re-colors the header row (bold black text on a blue fill instead of
# plain white-on-black).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the existing data row so the new "Abdul Quadir"
# record becomes row 2 and the current "Avinash" record shifts to row 3.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:G2").ClearFormats()

# --- Row 1: headers ---
# A1/B1 ("Corporate"/"employee") stay as-is; the rest change.
$ws.Range("C1").Value = "age"
$ws.Range("D1").Value = "sum insured"
$ws.Range("E1").Value = "A- AdditionD - DeletionC -Change"
$ws.Range("F1").Value = "Date of joining"
$ws.Range("G1").Value = "Date of Leaving"

# --- Row 2: new "Abdul Quadir" record ---
$ws.Range("A2").Value = "Demo Account"
$ws.Range("B2").Value = "Abdul Quadir"
$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 200000
$ws.Range("E2").Value = "D"
# Force text storage so these date-look-alike strings aren't silently
# converted into date serial numbers, then drop the number-format
# override again so the cell keeps plain (unstyled) formatting.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "01-01-1970"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "29-03-2023"
$ws.Range("G2").ClearFormats()

# --- Row 3: existing "Avinash" record, legend code updated ---
$ws.Range("E3").Value = "A"

# --- Header styling: bold black text on a blue fill (was plain white on
# black), applied across the now-wider A1:H1 header band. ---
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").Font.Color = 0
$ws.Range("A1:H1").Interior.Color = 16711680

# Match the saved selection/active cell (now the blank trailing column).
$ws.Range("H1").Select()
